# Add the new "summary_malfunction_reporting" field to the
# device_classification_fields sheet, as a new row inserted at row 6
# (pushing submission_type_id and everything below it down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device_classification_fields")

# Insert a new blank row at position 6; this shifts rows 6..18 down to
# 7..19 and copies formatting from the row above, matching how the
# other data rows in this sheet are styled.
$ws.Rows.Item(6).Insert()

# Populate the three cells of the new row.
$ws.Cells.Item(6, 2).Value = "summary_malfunction_reporting"
$ws.Cells.Item(6, 3).Value = "string"
$ws.Cells.Item(6, 4).Value = "The Voluntary Malfunction Summary Reporting Program allows participating companies to submit certain medical device malfunction reports in summary form on a quarterly basis.  The program applies to eligible devices regulated by the Center for Devices and Radiological Health (CDRH) and Center for Biologics Evaluation and Research (CBER), including device-led combination products.Value is one of the following:`nEligible = 510(K)`nIneligible = PMA"

# Match the row height used for the other short (non-wrapped-to-max)
# description rows.
$ws.Rows.Item(6).RowHeight = 102

# Widen column B slightly (description column lost its auto bestFit and
# became a manually sized column).
$ws.Columns.Item(2).ColumnWidth = 26.998697916666668

# Restore the view: zoomed in to 170%, with D6 (the new field's
# description cell) selected.
$excel.ActiveWindow.Zoom = 170
$ws.Range("D6").Select()
